$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '66.093.33'
$ws.Range("E2").Value = '  +6.88%  '
$ws.Range("D3").Value = '3.015.06'
$ws.Range("E3").Value = '  +4.35%  '
$ws.Range("E4").Value = '  +0.03%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '582.84'
$ws.Range("E5").Value = '  +2.75%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '161.27'
$ws.Range("E6").Value = '  +12.83%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.999'
$ws.Range("E7").Value = '  -0.17%  '
$ws.Range("D8").Value = '3.012.35'
$ws.Range("E8").Value = '  +4.37%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.518'
$ws.Range("E9").Value = '  +3.39%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '6.99'
$ws.Range("E10").Value = '  +1.41%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.155'
$ws.Range("E11").Value = '  +5.88%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.457'
$ws.Range("E12").Value = '  +6.21%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000251'
$ws.Range("E13").Value = '  +8.67%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '34.89'
$ws.Range("E14").Value = '  +9.51%  '
$ws.Range("E15").Value = '  +0.90%  '
$ws.Range("D16").Value = '66.087.44'
$ws.Range("E16").Value = '  +6.97%  '
$ws.Range("D17").Value = '3.516.74'
$ws.Range("E17").Value = '  +4.38%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '6.97'
$ws.Range("E18").Value = '  +6.82%  '
$ws.Range("D19").Value = '3.018.92'
$ws.Range("E19").Value = '  +3.63%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '459.36'
$ws.Range("E20").Value = '  +7.06%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13.94'
$ws.Range("E21").Value = '  +7.14%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.688'
$ws.Range("E22").Value = '  +5.49%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.34'
$ws.Range("E23").Value = '  +7.05%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '82.36'
$ws.Range("E24").Value = '  +4.33%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.28'
$ws.Range("E25").Value = '  +12.82%  '
$ws.Range("E26").Value = '  +3.08%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.67'
$ws.Range("E27").Value = '  +5.67%  '
$ws.Range("E28").Value = '  -0.09%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '8.20'
$ws.Range("E29").Value = '  +17.53%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.36'
$ws.Range("E30").Value = '  +17.30%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.0000105'
$ws.Range("E31").Value = '  -2.99%  '
$ws.Range("E32").Value = '  +4.17%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '27.01'
$ws.Range("E33").Value = '  +5.46%  '
$ws.Range("E34").Value = '  +3.01%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.998'
$ws.Range("E35").Value = '  -0.10%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.998'
$ws.Range("E36").Value = '  +5.33%  '
$ws.Range("E37").Value = '  +7.84%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.15'
$ws.Range("E38").Value = '  +13.01%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '50.11'
$ws.Range("E39").Value = '  +2.66%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.00'
$ws.Range("E40").Value = '  +7.44%  '
$ws.Range("B41").Value = 'TheGraph'
$ws.Range("C41").Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.305'
$ws.Range("E41").Value = '  +13.51%  '
$ws.Range("B42").Value = 'Kaspa'
$ws.Range("C42").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.122'
$ws.Range("E42").Value = '  +6.79%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '43.83'
$ws.Range("E43").Value = '  +10.13%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '8.52'
$ws.Range("E44").Value = '  +4.63%  '
$ws.Range("B45").Value = 'Bittensor'
$ws.Range("C45").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '384.86'
$ws.Range("E45").Value = '  +11.66%  '
$ws.Range("B46").Value = 'VeChain'
$ws.Range("C46").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0356'
$ws.Range("E46").Value = '  +5.74%  '
$ws.Range("B47").Value = 'Maker'
$ws.Range("C47").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D47").Value = '2.791.93'
$ws.Range("E47").Value = '  +3.93%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '134.83'
$ws.Range("E48").Value = '  +2.49%  '
$ws.Range("E49").Value = '  -0.06%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '23.91'
$ws.Range("E50").Value = '  +11.17%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.107'
$ws.Range("E51").Value = '  +4.35%  '
